$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking values stored as text. Force the
# cells to keep their text (string) storage before writing the new price so
# Excel doesn't silently re-type them as numbers.
$updates = @{
    2  = "266.06"
    3  = "22.62"
    4  = "6.205"
    6  = "3.561"
    7  = "6.712"
    9  = "0.8254"
    11 = "0.1593"
    12 = "0.08195"
    13 = "0.03396"
    14 = "0.03152"
    15 = "0.09237"
    16 = "3.909"
    17 = "0.001717"
    18 = "0.04799"
    19 = "0.006274"
    20 = "0.006276"
    22 = "0.0001500"
    23 = "3.709"
    24 = "2.260"
    27 = "0.0002681"
    40 = "0.04601"
    41 = "0.006981"
    42 = "0.1133"
    43 = "0.003130"
    44 = "0.01099"
    45 = "0.00006152"
    47 = "0.7700"
    48 = "0.2075"
    49 = "0.00002100"
    50 = "0.01240"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
